# TimeTrack.xlsx update: add new log entries (rows 5-11) and update totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time-tracking entries for rows 5 through 11, entered in the same order
# they were originally typed (rows 5-8 in order, then 10/9/11 interleaved).
$ws.Cells.Item(5, 1).Value = "Monday 16.4.18"
$ws.Cells.Item(5, 2).Value = "1630 - 1830"
$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 4).Value = "Rewriting assignment one in Java"

$ws.Cells.Item(6, 1).Value = "Tuesday 17.4.18"
$ws.Cells.Item(6, 2).Value = "1030 - 1130"
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = "Implementing LinkedList"

$ws.Cells.Item(7, 1).Value = "Monday 23.4.18"
$ws.Cells.Item(7, 2).Value = "1200-1500"
$ws.Cells.Item(7, 3).Value = 3
$ws.Cells.Item(7, 4).Value = "Theory + interface (ComparePolygons) + subclass (SortedPolygons)"

$ws.Cells.Item(8, 1).Value = "Monday 23.4.18"
$ws.Cells.Item(8, 2).Value = "1630-1730"
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = "Getting comparePolygons to work"

$ws.Cells.Item(10, 1).Value = "Sunday 29.4.18"
$ws.Cells.Item(9, 4).Value = "Insertion Sort"
$ws.Cells.Item(9, 1).Value = "Saturday 28.4.18"
$ws.Cells.Item(9, 2).Value = "2300-0100"
$ws.Cells.Item(9, 3).Value = 3

$ws.Cells.Item(10, 2).Value = "1200-1730"
$ws.Cells.Item(10, 3).Value = 5.5
$ws.Cells.Item(10, 4).Value = "Insertion Sort Works"

$ws.Cells.Item(11, 1).Value = "Sunday 29.4.1.8"
$ws.Cells.Item(11, 4).Value = "Comparable Interfaces"
$ws.Cells.Item(11, 3).Value = 2
$ws.Cells.Item(11, 2).Value = "2300-0100"

# Recalculate so the SUBTOTAL formula in C32 reflects the new hours.
$excel.Calculate()

# Update the selected cell to match the saved view state.
$ws.Range("B12").Select()
